$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.03940000000001
$ws.Range("E4").Value = 16.38939999999999
$ws.Range("D11").Value = -7.8339
$ws.Range("D12").Value = -6.9232
$ws.Range("E14").Value = 16.58090000000001
$ws.Range("D15").Value = -8.4922
$ws.Range("E26").Value = 16.38589999999998
$ws.Range("D27").Value = -8.441800000000002
$ws.Range("D28").Value = -8.347299999999997
$ws.Range("D31").Value = -7.7818
$ws.Range("E31").Value = 16.85560000000001
$ws.Range("D32").Value = -8.899999999999991
$ws.Range("E35").Value = 16.6496
$ws.Range("D36").Value = -8.177399999999999
$ws.Range("E37").Value = 16.5792
$ws.Range("D38").Value = -7.974000000000001
$ws.Range("E39").Value = 16.43350000000001
$ws.Range("E40").Value = 17.00460000000001
$ws.Range("E45").Value = 16.51609999999999
$ws.Range("D46").Value = -8.506099999999998
$ws.Range("E52").Value = 17.3686
$ws.Range("D54").Value = -8.213900000000006
$ws.Range("D55").Value = -8.367099999999999
$ws.Range("D56").Value = -8.047699999999995
$ws.Range("E57").Value = 16.66740000000001
$ws.Range("D67").Value = -6.337199999999996
$ws.Range("D69").Value = -6.909799999999997
$ws.Range("D72").Value = -7.422200000000001
$ws.Range("D73").Value = -7.605099999999994
$ws.Range("E81").Value = 15.93550000000001
$ws.Range("D83").Value = -8.676300000000005
$ws.Range("E83").Value = 16.57319999999999
$ws.Range("D86").Value = -7.391399999999993
$ws.Range("D91").Value = -6.564399999999999
$ws.Range("D93").Value = -6.561899999999998
$ws.Range("D99").Value = -7.660200000000006
$ws.Range("E100").Value = 16.4039
$ws.Range("E102").Value = 16.85029999999998
